# Boolean Do Suppliers Bid at Peak Capacity Factors.xlsx - US 3.3 commit files
$wb = $excel.ActiveWorkbook

# --- BDSBaPCF sheet edits ---
$ws = $wb.Worksheets.Item("BDSBaPCF")
$ws.Activate()

# Header label: clarify it is a boolean field; bold + wrap the header cell
$ws.Range("B1").Value = "Do Suppliers Bid at Peak Capacity Factors (Boolean)"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# petroleum / natural gas peaker no longer bid at peak capacity factors
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0

# "coal to gas" fuel row becomes "lignite" (still mirrors hard coal, row 2)
$ws.Range("A13").Value = "lignite"

# New fuel-type rows, each mirroring an existing boolean value
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"

# Selection left on B5 for this sheet
$ws.Range("B5").Select()

# --- Make the About sheet the active tab ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
